$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values for specific rows per the repull/mean recalculation
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("F8").Value = -9
$ws.Range("F9").Value = -4
$ws.Range("F10").Value = -5
$ws.Range("F11").Value = -5
$ws.Range("F12").Value = -7
